# "Se agrega alumno no valido" - add a (invalid/placeholder) new student
# row in the roster: name "Juan Perez " repeated in both the
# "Nombre y Apellido" and "nombre github" columns, classroom "Aula" = 0,
# and a mailto hyperlink in the "mail" column pointing at dire@gmail.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nombre y Apellido
$ws.Range("B10").Value = "Juan Perez "
# Aula
$ws.Range("C10").Value = 0
# mail - add hyperlink with visible text equal to the address
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:dire@gmail.com", [Type]::Missing, [Type]::Missing, "dire@gmail.com") | Out-Null
# Match the same "Hipervinculo" cell style used by the other mail cells
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
# nombre github
$ws.Range("E10").Value = "Juan Perez "

# Leave the active selection on the newly edited name cell
$ws.Range("B10").Select() | Out-Null
